$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3112
    5  = 2828
    6  = 190
    9  = 1515
    13 = 1258
    15 = 393
    17 = 53
    18 = 49
    21 = 95
    22 = 2801
    23 = 335
    24 = 8
    25 = 53
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
